# Add an "id" column as the new first column of the data sheet
# (Excel reports this sheet's tab as tabSelected/ActiveSheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing location name / latitude / longitude columns
# one place to the right, opening up a blank column A.
$ws.Columns("A:A").Insert()

# Give the new header cell (A1) the same formatting (border style)
# as the other header cells before filling in its text.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A1").Value = "id"

# Populate the id values for the three data rows.
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3

# Match the saved selection/active cell from the edit.
[void]$ws.Range("D4").Select()
